# Add two new columns ("I0" and "IF") to the sheet, mirroring the existing
# header/data layout (columns A-H) and extending it through column J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1) - same text style as the other headers (bold,
# bordered, centered). Easiest way to get an identical format is to copy
# the existing "IP" header cell's formatting onto the new header cells.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data cells (rows 2-3), plain numeric values with default formatting
# (same as the existing H2/H3 data cells).
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9
